$d = $word.ActiveDocument

$d.Content.Find.Execute("(123) 456-7890", $true, $false, $false, $false, $false, $true, 1, $false, "전화: (123) 456-7890", 2)

$d.Content.Find.Execute("수석 애니메이터(2018년 1월~현재)", $true, $false, $false, $false, $false, $true, 1, $false, "ABC 스튜디오: 리드 애니메이터(2018년 1월 - 현재)", 2)

$d.Content.Find.Execute("선임 애니메이터(2015년 1월~2017년 12월)", $true, $false, $false, $false, $false, $true, 1, $false, "XYZ 미디어: 선임 애니메이터(2015년 6월 - 2017년 12월)", 2)

$d.Content.Find.Execute("보조 애니메이터(2012년 9월~2015년 5월)", $true, $false, $false, $false, $false, $true, 1, $false, "MNO 엔터테인먼트: 주니어 애니메이터 (2012년 9월 - 2015년 5월)", 2)
